# Applies two changes described by the commit diff:
#   1. The table on slide 16 switches from the custom "Table_0" style
#      ({BFB9E666-8343-4D68-9134-4E6F5BDF7B72}) to the built-in table
#      style {32F3750B-7429-42F8-B5A1-E73CB89BCADC}.
#   2. The presentation's theme (ppt/theme/theme1.xml, the theme behind
#      the slide master / "Integral" design) is replaced by what used to
#      be the secondary "Office Theme" palette - i.e. the clrScheme
#      colours swap from the Integral palette to the stock Office palette
#      (font scheme / format scheme are already identical between the two
#      themes, only the colour slots differ).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$s = $p.Slides.Item(16)
$shp = $s.Shapes.Item(3)
$tbl = $shp.Table
$tbl.ApplyStyle("{32F3750B-7429-42F8-B5A1-E73CB89BCADC}")

# --- 2. Theme colours -------------------------------------------------
# ThemeColorScheme.Colors index order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
# 11 hlink, 12 folHlink
$officeColors = @{
    3  = 0x44546A  # dk2
    4  = 0xE7E6E6  # lt2
    5  = 0x5B9BD5  # accent1
    6  = 0xED7D31  # accent2
    7  = 0xA5A5A5  # accent3
    8  = 0xFFC000  # accent4
    9  = 0x4472C4  # accent5
    10 = 0x70AD47  # accent6
    11 = 0x0563C1  # hlink
    12 = 0x954F72  # folHlink
}

$ccs = $p.SlideMaster.Theme.ThemeColorScheme
foreach ($idx in $officeColors.Keys) {
    # $hex is a 0xRRGGBB literal - pull bytes out in that order, then
    # repack as the OLE RGB() integer (R | G<<8 | B<<16) the host expects.
    $hex = $officeColors[$idx]
    $r = ($hex -shr 16) -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $b = $hex -band 0xFF
    $rgb = $r -bor ($g -shl 8) -bor ($b -shl 16)
    $ccs.Colors($idx).RGB = $rgb
}
